$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the text values of A7 and A8 (Climate <-> Environmental)
$a7 = $ws.Range("A7").Text
$a8 = $ws.Range("A8").Text
$ws.Range("A7").Value = $a8
$ws.Range("A8").Value = $a7

# Update the selection to match the new active range (A7:A8, active cell A7)
$ws.Range("A7:A8").Select()
